$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 (UCL) ROI values refreshed as part of the 04.02.2025 update.
# Cells whose value crosses zero also flip their red-fill "positive ROI"
# highlight (style index 2 in the workbook == Interior.ColorIndex 3),
# while cells settling back to a non-positive value drop back to Normal.

$ws.Range("B30").Value = 32.88

$ws.Range("D30").Style = "Normal"
$ws.Range("D30").Value = -21.04

$ws.Range("E30").Value = -0.5

$ws.Range("F30").Value = 5.4

$ws.Range("G30").Value = -36.36

$ws.Range("H30").Value = -9.2

$ws.Range("I30").Value = 6

$ws.Range("J30").Interior.ColorIndex = 3
$ws.Range("J30").Value = 1.07

$ws.Range("K30").Value = -27.59

$ws.Range("L30").Value = -28.56

$ws.Range("M30").Value = -10.51

$ws.Range("N30").Value = 27.55

$ws.Range("O30").Value = -11.52

$ws.Range("P30").Value = 13.35

$ws.Range("Q30").Value = 8.9

$ws.Range("R30").Value = 9.24

$ws.Range("S30").Style = "Normal"
$ws.Range("S30").Value = -6.88

$ws.Range("T30").Interior.ColorIndex = 3
$ws.Range("T30").Value = 11.69

$ws.Range("U30").Style = "Normal"
$ws.Range("U30").Value = -7.37

$ws.Range("V30").Value = -27.25

$ws.Range("W30").Style = "Normal"
$ws.Range("W30").Value = -14.11

$ws.Range("X30").Interior.ColorIndex = 3
$ws.Range("X30").Value = 3.71

$ws.Range("Y30").Style = "Normal"
$ws.Range("Y30").Value = -3.62

$ws.Range("Z30").Value = -32

$ws.Range("AA30").Value = -8.17

$ws.Range("AB30").Value = 51
